# Auto-generated edit script applying the crawl-refresh diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking id/price columns stay stored as text (matches source data)
# via Range.NumberFormat = "@" immediately before assigning the value.

$ws.Cells.Item(2, 15).Value2 = '2022-09-14 21:00:47'

$ws.Range("A3").NumberFormat = "@"
$ws.Range("H3").NumberFormat = "@"
$ws.Cells.Item(3, 1).Value2 = '6973029'
$ws.Cells.Item(3, 2).Value2 = 'Severin Tischgrill PG 8565'
$ws.Cells.Item(3, 3).Value2 = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/severin-tischgrill-pg-8565/p/6973029'
$ws.Cells.Item(3, 5).Value2 = 1
$ws.Cells.Item(3, 6).Value2 = 3
$ws.Cells.Item(3, 7).Value2 = 'Severin'
$ws.Cells.Item(3, 8).Value2 = '49.95'
$ws.Cells.Item(3, 13).Value2 = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''kuechengeraete'']'
$ws.Cells.Item(3, 14).Value2 = 'Severin Tischgrill PG 8565 50% Aktion 49.95 Schweizer Franken statt 99.95 Schweizer Franken'
$ws.Cells.Item(3, 15).Value2 = '2022-09-14 21:00:47'

$ws.Cells.Item(4, 15).Value2 = '2022-09-14 21:00:47'

$ws.Cells.Item(5, 15).Value2 = '2022-09-14 21:00:47'

$ws.Cells.Item(6, 15).Value2 = '2022-09-14 21:00:47'

$ws.Cells.Item(7, 15).Value2 = '2022-09-14 21:00:47'

$ws.Cells.Item(8, 15).Value2 = '2022-09-14 21:00:47'

$ws.Cells.Item(9, 15).Value2 = '2022-09-14 21:00:47'

$ws.Cells.Item(10, 15).Value2 = '2022-09-14 21:00:47'

$ws.Cells.Item(11, 15).Value2 = '2022-09-14 21:00:47'

$ws.Cells.Item(12, 15).Value2 = '2022-09-14 21:00:47'

$ws.Cells.Item(13, 15).Value2 = '2022-09-14 21:00:47'

$ws.Cells.Item(14, 15).Value2 = '2022-09-14 21:00:47'

$ws.Range("A15").NumberFormat = "@"
$ws.Range("H15").NumberFormat = "@"
$ws.Cells.Item(15, 1).Value2 = '3862219'
$ws.Cells.Item(15, 2).Value2 = 'Varta Active LED Camping Lanterne'
$ws.Cells.Item(15, 3).Value2 = '/de/haushalt-tier/elektroartikel-batterien/beleuchtung/taschenlampen-lichter/varta-active-led-camping-lanterne/p/3862219'
$ws.Cells.Item(15, 7).Value2 = 'Varta'
$ws.Cells.Item(15, 8).Value2 = '38.50'
$ws.Cells.Item(15, 13).Value2 = '[''haushalt-tier'', ''elektroartikel-batterien'', ''beleuchtung'', ''taschenlampen-lichter'']'
$ws.Cells.Item(15, 14).Value2 = 'Varta Active LED Camping Lanterne 38.50 Schweizer Franken'
$ws.Cells.Item(15, 15).Value2 = '2022-09-14 21:00:47'

$ws.Range("A16").NumberFormat = "@"
$ws.Range("H16").NumberFormat = "@"
$ws.Cells.Item(16, 1).Value2 = '6725106'
$ws.Cells.Item(16, 2).Value2 = 'Braun Dampfbügelstation IS1012BL'
$ws.Cells.Item(16, 3).Value2 = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/braun-dampfbuegelstation-is1012bl/p/6725106'
$ws.Cells.Item(16, 7).Value2 = 'Braun'
$ws.Cells.Item(16, 8).Value2 = '149.00'
$ws.Cells.Item(16, 13).Value2 = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''haushaltsgeraete-kabel'']'
$ws.Cells.Item(16, 14).Value2 = 'Braun Dampfbügelstation IS1012BL 149.00 Schweizer Franken'
$ws.Cells.Item(16, 15).Value2 = '2022-09-14 21:00:47'

$ws.Range("A17").NumberFormat = "@"
$ws.Range("H17").NumberFormat = "@"
$ws.Cells.Item(17, 1).Value2 = '6995204'
$ws.Cells.Item(17, 2).Value2 = 'Electrolux AirFryer Range Explore 6'
$ws.Cells.Item(17, 3).Value2 = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/electrolux-airfryer-range-explore-6/p/6995204'
$ws.Cells.Item(17, 7).Value2 = 'Electrolux'
$ws.Cells.Item(17, 8).Value2 = '99.50'
$ws.Cells.Item(17, 14).Value2 = 'Electrolux AirFryer Range Explore 6 50% Aktion 99.50 Schweizer Franken statt 199.00 Schweizer Franken'
$ws.Cells.Item(17, 15).Value2 = '2022-09-14 21:00:47'

$ws.Range("A18").NumberFormat = "@"
$ws.Range("H18").NumberFormat = "@"
$ws.Cells.Item(18, 1).Value2 = '6982305'
$ws.Cells.Item(18, 2).Value2 = 'Kenwood Küchenmaschine MultiOne KHH326WH'
$ws.Cells.Item(18, 3).Value2 = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/kenwood-kuechenmaschine-multione-khh326wh/p/6982305'
$ws.Cells.Item(18, 5).Value2 = ''
$ws.Cells.Item(18, 6).Value2 = 0
$ws.Cells.Item(18, 7).Value2 = 'Kenwood'
$ws.Cells.Item(18, 8).Value2 = '299.00'
$ws.Cells.Item(18, 13).Value2 = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''kuechengeraete'']'
$ws.Cells.Item(18, 14).Value2 = 'Kenwood Küchenmaschine MultiOne KHH326WH 40% Aktion 299.00 Schweizer Franken statt 499.00 Schweizer Franken'
$ws.Cells.Item(18, 15).Value2 = '2022-09-14 21:00:47'

$ws.Range("A19").NumberFormat = "@"
$ws.Range("H19").NumberFormat = "@"
$ws.Cells.Item(19, 1).Value2 = '6735643'
$ws.Cells.Item(19, 2).Value2 = 'LED 31V Anschlussset Transf.+Verl.kabel'
$ws.Cells.Item(19, 3).Value2 = '/de/haushalt-tier/haushalt-kueche/uebrige-haushaltsartikel/led-31v-anschlussset-transfverlkabel/p/6735643'
$ws.Cells.Item(19, 5).Value2 = 1
$ws.Cells.Item(19, 6).Value2 = 5
$ws.Cells.Item(19, 7).Value2 = 'Coop'
$ws.Cells.Item(19, 8).Value2 = '9.95'
$ws.Cells.Item(19, 13).Value2 = '[''haushalt-tier'', ''haushalt-kueche'', ''uebrige-haushaltsartikel'']'
$ws.Cells.Item(19, 14).Value2 = 'LED 31V Anschlussset Transf.+Verl.kabel 50% Aktion 9.95 Schweizer Franken statt 19.95 Schweizer Franken'
$ws.Cells.Item(19, 15).Value2 = '2022-09-14 21:00:47'

$ws.Range("A20").NumberFormat = "@"
$ws.Range("H20").NumberFormat = "@"
$ws.Cells.Item(20, 1).Value2 = '5882124'
$ws.Cells.Item(20, 2).Value2 = 'Philips Avent Audio Monitors DECT-Babyphone'
$ws.Cells.Item(20, 3).Value2 = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/philips-avent-audio-monitors-dect-babyphone/p/5882124'
$ws.Cells.Item(20, 7).Value2 = 'Avent'
$ws.Cells.Item(20, 8).Value2 = '99.90'
$ws.Cells.Item(20, 13).Value2 = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'']'
$ws.Cells.Item(20, 14).Value2 = 'Philips Avent Audio Monitors DECT-Babyphone 99.90 Schweizer Franken'
$ws.Cells.Item(20, 15).Value2 = '2022-09-14 21:00:47'

$ws.Cells.Item(21, 15).Value2 = '2022-09-14 21:00:47'

$ws.Cells.Item(22, 15).Value2 = '2022-09-14 21:00:47'

$ws.Cells.Item(23, 15).Value2 = '2022-09-14 21:00:47'

$ws.Cells.Item(24, 15).Value2 = '2022-09-14 21:00:47'

$ws.Cells.Item(25, 15).Value2 = '2022-09-14 21:00:47'

$ws.Cells.Item(26, 15).Value2 = '2022-09-14 21:00:47'

$ws.Cells.Item(27, 15).Value2 = '2022-09-14 21:00:47'

$ws.Cells.Item(28, 15).Value2 = '2022-09-14 21:00:47'

$ws.Cells.Item(29, 15).Value2 = '2022-09-14 21:00:47'

$ws.Cells.Item(30, 15).Value2 = '2022-09-14 21:00:47'

$ws.Cells.Item(31, 15).Value2 = '2022-09-14 21:00:47'

$ws.Cells.Item(32, 15).Value2 = '2022-09-14 21:00:47'

$ws.Cells.Item(33, 15).Value2 = '2022-09-14 21:00:47'

$ws.Range("A34").NumberFormat = "@"
$ws.Range("H34").NumberFormat = "@"
$ws.Range("K34").NumberFormat = "@"
$ws.Cells.Item(34, 1).Value2 = '4589933'
$ws.Cells.Item(34, 2).Value2 = 'Varta Longlife AAA 4er Bli'
$ws.Cells.Item(34, 3).Value2 = '/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/varta-longlife-aaa-4er-bli/p/4589933'
$ws.Cells.Item(34, 4).Value2 = '4ST'
$ws.Cells.Item(34, 8).Value2 = '7.95'
$ws.Cells.Item(34, 9).Value2 = '1.99/1ST'
$ws.Cells.Item(34, 10).Value2 = 'Preis pro 1 Stück'
$ws.Cells.Item(34, 11).Value2 = '1.99'
$ws.Cells.Item(34, 12).Value2 = '1ST'
$ws.Cells.Item(34, 13).Value2 = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''aaa'']'
$ws.Cells.Item(34, 14).Value2 = 'Varta Longlife AAA 4er Bli 7.95 Schweizer Franken'
$ws.Cells.Item(34, 15).Value2 = '2022-09-14 21:00:47'

$ws.Range("A35").NumberFormat = "@"
$ws.Range("H35").NumberFormat = "@"
$ws.Range("K35").NumberFormat = "@"
$ws.Cells.Item(35, 1).Value2 = '4589935'
$ws.Cells.Item(35, 2).Value2 = 'Varta Longlife C 2er Bli'
$ws.Cells.Item(35, 3).Value2 = '/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-longlife-c-2er-bli/p/4589935'
$ws.Cells.Item(35, 4).Value2 = '2ST'
$ws.Cells.Item(35, 8).Value2 = '6.95'
$ws.Cells.Item(35, 9).Value2 = '3.48/1ST'
$ws.Cells.Item(35, 11).Value2 = '3.48'
$ws.Cells.Item(35, 13).Value2 = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''andere-batterien'']'
$ws.Cells.Item(35, 14).Value2 = 'Varta Longlife C 2er Bli 6.95 Schweizer Franken'
$ws.Cells.Item(35, 15).Value2 = '2022-09-14 21:00:47'

$ws.Range("A36").NumberFormat = "@"
$ws.Range("H36").NumberFormat = "@"
$ws.Range("K36").NumberFormat = "@"
$ws.Cells.Item(36, 1).Value2 = '3494067'
$ws.Cells.Item(36, 2).Value2 = 'Varta Longlife Max Power AAA 4er Bli'
$ws.Cells.Item(36, 3).Value2 = '/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/varta-longlife-max-power-aaa-4er-bli/p/3494067'
$ws.Cells.Item(36, 4).Value2 = '4ST'
$ws.Cells.Item(36, 8).Value2 = '9.95'
$ws.Cells.Item(36, 9).Value2 = '2.49/1ST'
$ws.Cells.Item(36, 11).Value2 = '2.49'
$ws.Cells.Item(36, 13).Value2 = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''aaa'']'
$ws.Cells.Item(36, 14).Value2 = 'Varta Longlife Max Power AAA 4er Bli 9.95 Schweizer Franken'
$ws.Cells.Item(36, 15).Value2 = '2022-09-14 21:00:47'

# Row 37 (old trailing product) no longer exists; remaining rows shifted up by one
$ws.Rows.Item(37).Delete()

